$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.648.70"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").Value = "1.594.53"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'211.07"

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").Value = "1.819.14"
$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").Value = "1.584.55"
$ws.Range("E13").Value = "  -2.52%  "

$ws.Range("E14").Value = "  -2.33%  "

$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -3.04%  "

$ws.Range("D16").Value = "'64.76"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "26.632.73"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").Value = "'209.01"
$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = "  -2.19%  "

$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").Value = "'146.63"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = "  -4.14%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "'15.30"
$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").Value = "'3.21"
$ws.Range("E32").Value = "  -2.83%  "

$ws.Range("E33").Value = "  -4.32%  "

$ws.Range("E34").Value = "  -2.97%  "

$ws.Range("D35").Value = "1.292.93"
$ws.Range("E35").Value = "  -3.13%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("E37").Value = "  -4.90%  "

$ws.Range("E38").Value = "  -2.66%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("D41").Value = "'0.791"
$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -1.65%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "'63.51"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").Value = "1.731.36"
$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("D46").Value = "'0.895"
$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("D47").Value = "'89.69"
$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("D48").Value = "'1.63"
$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").Value = "'0.0981"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D51").Value = "'7.48"
$ws.Range("E51").Value = "  -1.06%  "
